# Regenerate orders with updated distance/size labels.
# The experiment's viewing distances and one stimulus size code were
# renamed; every cell that encodes them (Condition, Filename_Left,
# Filename_Right, Distance, Size - and any filenames/labels built from
# them) needs the same substring substitution applied:
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Distances (mm) were bumped up:
$used.Replace("D80", "D86")
$used.Replace("D64", "D69")
$used.Replace("D51", "D55")

# One of the size codes changed too:
$used.Replace("S30", "S31")
